# Update with restock suggestion
#
# - "Forecast Comparison" sheet:
#     * Fills in the Week_Start_Date (col B) for every data row.
#     * Refreshes Inventory Coverage (L) and Seasonality Index (P).
#     * Rows W8..W16 flip Stockout Risk (M) / Reorder Urgency (N)
#       from Low/Normal to High/Urgent (inventory has run out).
#     * Drops the "Sales Volume Rank" column (Q) entirely, and renames
#       the old "Lifecycle Stage" column (R) to take its place as the
#       new Q, with every row's stage recomputed to "Decline".
# - "Summary" sheet:
#     * Max/Min Forecast Week become "N/A" (no longer meaningful once
#       the product is in decline / restock territory).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")

# Per-row data: row number, Week_Start_Date, Inventory Coverage,
# Stockout Risk, Reorder Urgency, Seasonality Index.
$rows = @(
    @{ Row = 2;  Date = "2025-02-02"; Coverage = 8;    Risk = "Low";  Urgency = "Normal"; Seasonality = 1.04 },
    @{ Row = 3;  Date = "2025-02-09"; Coverage = 7;    Risk = "Low";  Urgency = "Normal"; Seasonality = 0.88 },
    @{ Row = 4;  Date = "2025-02-16"; Coverage = 6;    Risk = "Low";  Urgency = "Normal"; Seasonality = 0.91 },
    @{ Row = 5;  Date = "2025-02-23"; Coverage = 4.17; Risk = "Low";  Urgency = "Normal"; Seasonality = 0.95 },
    @{ Row = 6;  Date = "2025-03-02"; Coverage = 3.17; Risk = "Low";  Urgency = "Normal"; Seasonality = 0.83 },
    @{ Row = 7;  Date = "2025-03-09"; Coverage = 2.17; Risk = "Low";  Urgency = "Normal"; Seasonality = 0.97 },
    @{ Row = 8;  Date = "2025-03-16"; Coverage = 1.17; Risk = "Low";  Urgency = "Normal"; Seasonality = 1.12 },
    @{ Row = 9;  Date = "2025-03-23"; Coverage = 0.17; Risk = "High"; Urgency = "Urgent"; Seasonality = 0.92 },
    @{ Row = 10; Date = "2025-03-30"; Coverage = 0;    Risk = "High"; Urgency = "Urgent"; Seasonality = 1.06 },
    @{ Row = 11; Date = "2025-04-06"; Coverage = 0;    Risk = "High"; Urgency = "Urgent"; Seasonality = 0.96 },
    @{ Row = 12; Date = "2025-04-13"; Coverage = 0;    Risk = "High"; Urgency = "Urgent"; Seasonality = 1.11 },
    @{ Row = 13; Date = "2025-04-20"; Coverage = 0;    Risk = "High"; Urgency = "Urgent"; Seasonality = 1.01 },
    @{ Row = 14; Date = "2025-04-27"; Coverage = 0;    Risk = "High"; Urgency = "Urgent"; Seasonality = 1.16 },
    @{ Row = 15; Date = "2025-05-04"; Coverage = 0;    Risk = "High"; Urgency = "Urgent"; Seasonality = 0.85 },
    @{ Row = 16; Date = "2025-05-11"; Coverage = 0;    Risk = "High"; Urgency = "Urgent"; Seasonality = 1.04 },
    @{ Row = 17; Date = "2025-05-18"; Coverage = 0;    Risk = "High"; Urgency = "Urgent"; Seasonality = 0.93 }
)

foreach ($r in $rows) {
    $rowNum = $r.Row

    # Column B - Week_Start_Date. Force text storage (leading apostrophe)
    # so "2025-02-02" isn't auto-converted into a date serial.
    $ws.Cells.Item($rowNum, 2).Value = "'" + $r.Date

    # Column L - Inventory Coverage
    $ws.Cells.Item($rowNum, 12).Value = $r.Coverage

    # Column M - Stockout Risk
    $ws.Cells.Item($rowNum, 13).Value = $r.Risk

    # Column N - Reorder Urgency
    $ws.Cells.Item($rowNum, 14).Value = $r.Urgency

    # Column P - Seasonality Index
    $ws.Cells.Item($rowNum, 16).Value = $r.Seasonality
}

# Drop column Q ("Sales Volume Rank"); this shifts the old column R
# ("Lifecycle Stage") left into Q, and the dimension shrinks to A1:Q17.
$ws.Columns.Item(17).Delete()

# Recompute the (now-Q) Lifecycle Stage for every row to "Decline".
$ws.Cells.Item(1, 17).Value = "Lifecycle Stage"
for ($rowNum = 2; $rowNum -le 17; $rowNum++) {
    $ws.Cells.Item($rowNum, 17).Value = "Decline"
}

# Summary sheet: the max/min forecast week no longer apply.
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B13").Value = "N/A"
$summary.Range("B15").Value = "N/A"
